# This script reproduces a new weekly price record for "Albahaca" at
# "Terminal La Palmera de La Serena" being inserted into the data table.
# A new row is inserted at row 127 (pushing the former rows 127-163 down
# to 128-164) and populated with the new observation's values, while all
# the other columns (A,B,C,E,F,G,H,I,R) keep the same constants used
# throughout this data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 127, shifting existing rows down.
$ws.Rows.Item(127).Insert()

# Fill in the new row 127 with the new data point.
$ws.Range("A127").Value = 8
$ws.Range("B127").Value = "Terminal La Palmera de La Serena"
$ws.Range("C127").Value = "Coquimbo"
$ws.Range("D127").Value = (Get-Date -Year 2023 -Month 3 -Day 3 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E127").Value = 4
$ws.Range("F127").Value = 100112052
$ws.Range("G127").Value = "Albahaca"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 2800
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 2900
$ws.Range("N127").Value = "`$/docena de matas"
$ws.Range("O127").Value = "Provincia del Elquí"
$ws.Range("P127").Value = 483
$ws.Range("Q127").Value = 6
$ws.Range("R127").Value = "Hortaliza"
